$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column G (7th column, "I/O") so it matches the width used by
# columns F ("DSP") and H ("XADC") -- stored width 10.3125
$ws.Columns.Item(7).ColumnWidth = 9.5

# Updated power figures for row 2 (added 64 bit designs)
$ws.Range("B2").Value = 0.02843075431883335
$ws.Range("C2").Value = 0.0077955471351742744
$ws.Range("D2").Value = 0.004497642163187265
$ws.Range("E2").Value = 0.002616662997752428
$ws.Range("F2").Value = 0.0000000353668170305354579
$ws.Range("G2").Value = 0.0007905749953351915
$ws.Range("I2").Value = 1.2575732469558716
$ws.Range("J2").Value = 0.12738706171512604
$ws.Range("K2").Value = 1.4296834468841553
